# Add more draft picks (fans/thebat depth-chart additions) to the
# "draftpicks" sheet: rows 301-308.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("draftpicks")

# team, player, salary, position, drafted (Excel date serial)
$rows = @(
    @("hermanos", "Rougned Odor",     19, "MI", 43124),
    @("bears",    "Andrew McCutchen", 23, "OF", 43124),
    @("pkdodgers","Alex Verdugo",      8, "OF", 43124),
    @("deano",    "Avisail Garcia",    8, "OF", 43124),
    @("ottawa",   "Hunter Renfroe",   14, "OF", 43124),
    @("chicago",  "Tyler Flowers",     3, "C",  43124),
    @("deano",    "Corey Knebel",     16, "P",  43124),
    @("ottawa",   "Mike Zunino",       8, "C",  43124)
)

$startRow = 301
# Reference cells whose formatting we reuse so no new cell styles get
# fabricated: B293 carries the "no wrap" player-name style used for the
# new rows, and E300 carries the short-date format used throughout
# column E.
$playerFmtSource = $ws.Cells.Item(293, 2)
$dateFmtSource   = $ws.Cells.Item(300, 5)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $team     = $rows[$i][0]
    $player   = $rows[$i][1]
    $salary   = $rows[$i][2]
    $position = $rows[$i][3]
    $drafted  = $rows[$i][4]

    $ws.Cells.Item($r, 1).Value = $team

    $playerFmtSource.Copy()
    $ws.Cells.Item($r, 2).PasteSpecial(-4122) # xlPasteFormats
    $ws.Cells.Item($r, 2).Value = $player

    $ws.Cells.Item($r, 3).Value = $salary

    $ws.Cells.Item($r, 4).Value = $position

    $dateFmtSource.Copy()
    $ws.Cells.Item($r, 5).PasteSpecial(-4122) # xlPasteFormats
    $ws.Cells.Item($r, 5).Value = $drafted
}

$excel.CutCopyMode = $false

# Update the view to mirror the author's workbook state after the edit.
$ws.Range("I295").Select()
